$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 38740
$ws.Cells.Item(2, 5).Value = 5638
$ws.Cells.Item(2, 6).Value = 5638
$ws.Cells.Item(2, 7).Value = 5302
$ws.Cells.Item(2, 8).Value = 3851
$ws.Cells.Item(2, 9).Value = 3791
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 38546
$ws.Cells.Item(2, 12).Value = 9763
$ws.Cells.Item(2, 13).Value = 28783
$ws.Cells.Item(2, 14).Value = 28587
$ws.Cells.Item(2, 15).Value = 196
$ws.Cells.Item(2, 16).Value = 345
$ws.Cells.Item(2, 17).Value = 5861
$ws.Cells.Item(2, 18).Value = -4907
$ws.Cells.Item(2, 19).Value = -463
$ws.Cells.Item(2, 20).Value = 2672
$ws.Cells.Item(2, 21).Value = 3189
$ws.Cells.Item(2, 22).Value = 1544
$ws.Cells.Item(2, 23).Value = 14.55
$ws.Cells.Item(2, 24).Value = 9.94
$ws.Cells.Item(2, 25).Value = 14.01
$ws.Cells.Item(2, 26).Value = 10.62
$ws.Cells.Item(2, 27).Value = 33.92
$ws.Cells.Item(2, 28).Value = 8278.780000000001
$ws.Cells.Item(2, 29).Value = 5493
$ws.Cells.Item(2, 30).Value = 40.42
$ws.Cells.Item(2, 31).Value = 41444
$ws.Cells.Item(2, 32).Value = 5.36
$ws.Cells.Item(2, 33).Value = 900
$ws.Cells.Item(2, 34).Value = 0.41
$ws.Cells.Item(2, 35).Value = 16.39
$ws.Cells.Item(2, 36).Value = 58458490

# Row 3
$ws.Cells.Item(3, 4).Value = 47666
$ws.Cells.Item(3, 5).Value = 7729
$ws.Cells.Item(3, 6).Value = 7729
$ws.Cells.Item(3, 7).Value = 7783
$ws.Cells.Item(3, 8).Value = 5848
$ws.Cells.Item(3, 9).Value = 5775
$ws.Cells.Item(3, 10).Value = 73
$ws.Cells.Item(3, 11).Value = 44431
$ws.Cells.Item(3, 12).Value = 10708
$ws.Cells.Item(3, 13).Value = 33723
$ws.Cells.Item(3, 14).Value = 33456
$ws.Cells.Item(3, 15).Value = 267
$ws.Cells.Item(3, 16).Value = 345
$ws.Cells.Item(3, 17).Value = 6275
$ws.Cells.Item(3, 18).Value = -2169
$ws.Cells.Item(3, 19).Value = -694
$ws.Cells.Item(3, 20).Value = 2153
$ws.Cells.Item(3, 21).Value = 4122
$ws.Cells.Item(3, 22).Value = 1516
$ws.Cells.Item(3, 23).Value = 16.21
$ws.Cells.Item(3, 24).Value = 12.27
$ws.Cells.Item(3, 25).Value = 18.62
$ws.Cells.Item(3, 26).Value = 14.09
$ws.Cells.Item(3, 27).Value = 31.75
$ws.Cells.Item(3, 28).Value = 9689.139999999999
$ws.Cells.Item(3, 29).Value = 8368
$ws.Cells.Item(3, 30).Value = 49.53
$ws.Cells.Item(3, 31).Value = 48503
$ws.Cells.Item(3, 32).Value = 8.550000000000001
$ws.Cells.Item(3, 33).Value = 1350
$ws.Cells.Item(3, 34).Value = 0.33
$ws.Cells.Item(3, 35).Value = 16.13
$ws.Cells.Item(3, 36).Value = 58458490

# Row 4
$ws.Cells.Item(4, 4).Value = 56454
$ws.Cells.Item(4, 5).Value = 8481
$ws.Cells.Item(4, 6).Value = 8481
$ws.Cells.Item(4, 7).Value = 8566
$ws.Cells.Item(4, 8).Value = 6457
$ws.Cells.Item(4, 9).Value = 6393
$ws.Cells.Item(4, 10).Value = 65
$ws.Cells.Item(4, 11).Value = 51816
$ws.Cells.Item(4, 12).Value = 12849
$ws.Cells.Item(4, 13).Value = 38966
$ws.Cells.Item(4, 14).Value = 38788
$ws.Cells.Item(4, 15).Value = 178
$ws.Cells.Item(4, 16).Value = 345
$ws.Cells.Item(4, 17).Value = 6770
$ws.Cells.Item(4, 18).Value = -6096
$ws.Cells.Item(4, 19).Value = -1330
$ws.Cells.Item(4, 20).Value = 4948
$ws.Cells.Item(4, 21).Value = 1821
$ws.Cells.Item(4, 22).Value = 1477
$ws.Cells.Item(4, 23).Value = 15.02
$ws.Cells.Item(4, 24).Value = 11.44
$ws.Cells.Item(4, 25).Value = 17.7
$ws.Cells.Item(4, 26).Value = 13.42
$ws.Cells.Item(4, 27).Value = 32.98
$ws.Cells.Item(4, 28).Value = 11241.22
$ws.Cells.Item(4, 29).Value = 9262
$ws.Cells.Item(4, 30).Value = 34.71
$ws.Cells.Item(4, 31).Value = 56234
$ws.Cells.Item(4, 32).Value = 5.72
$ws.Cells.Item(4, 33).Value = 1580
$ws.Cells.Item(4, 34).Value = 0.49
$ws.Cells.Item(4, 35).Value = 17.06
$ws.Cells.Item(4, 36).Value = 58458490

# Row 5
$ws.Cells.Item(5, 4).Value = 51238
$ws.Cells.Item(5, 5).Value = 5964
$ws.Cells.Item(5, 6).Value = 5964
$ws.Cells.Item(5, 7).Value = 5673
$ws.Cells.Item(5, 8).Value = 3980
$ws.Cells.Item(5, 9).Value = 3940
$ws.Cells.Item(5, 10).Value = 40
$ws.Cells.Item(5, 11).Value = 53757
$ws.Cells.Item(5, 12).Value = 12003
$ws.Cells.Item(5, 13).Value = 41754
$ws.Cells.Item(5, 14).Value = 41538
$ws.Cells.Item(5, 15).Value = 216
$ws.Cells.Item(5, 16).Value = 345
$ws.Cells.Item(5, 17).Value = 5126
$ws.Cells.Item(5, 18).Value = -4008
$ws.Cells.Item(5, 19).Value = -1109
$ws.Cells.Item(5, 20).Value = 7686
$ws.Cells.Item(5, 21).Value = -2560
$ws.Cells.Item(5, 22).Value = 1295
$ws.Cells.Item(5, 23).Value = 11.64
$ws.Cells.Item(5, 24).Value = 7.77
$ws.Cells.Item(5, 25).Value = 9.81
$ws.Cells.Item(5, 26).Value = 7.54
$ws.Cells.Item(5, 27).Value = 28.75
$ws.Cells.Item(5, 28).Value = 12104.31
$ws.Cells.Item(5, 29).Value = 5709
$ws.Cells.Item(5, 30).Value = 53.34
$ws.Cells.Item(5, 31).Value = 60220
$ws.Cells.Item(5, 32).Value = 5.06
$ws.Cells.Item(5, 33).Value = 1280
$ws.Cells.Item(5, 34).Value = 0.42
$ws.Cells.Item(5, 35).Value = 22.42
$ws.Cells.Item(5, 36).Value = 58458490

# Row 6
$ws.Cells.Item(6, 4).Value = 52778
$ws.Cells.Item(6, 5).Value = 4820
$ws.Cells.Item(6, 6).Value = 4820
$ws.Cells.Item(6, 7).Value = 4530
$ws.Cells.Item(6, 8).Value = 3348
$ws.Cells.Item(6, 9).Value = 3322
$ws.Cells.Item(6, 11).Value = 53711
$ws.Cells.Item(6, 12).Value = 9247
$ws.Cells.Item(6, 13).Value = 44464
$ws.Cells.Item(6, 14).Value = 44220
$ws.Cells.Item(6, 16).Value = 345
$ws.Cells.Item(6, 17).Value = 6467
$ws.Cells.Item(6, 18).Value = -4139
$ws.Cells.Item(6, 19).Value = -1034
$ws.Cells.Item(6, 20).Value = 4055
$ws.Cells.Item(6, 21).Value = 2412
$ws.Cells.Item(6, 22).Value = 1161
$ws.Cells.Item(6, 23).Value = 9.130000000000001
$ws.Cells.Item(6, 24).Value = 6.34
$ws.Cells.Item(6, 25).Value = 7.75
$ws.Cells.Item(6, 26).Value = 6.23
$ws.Cells.Item(6, 27).Value = 20.8
$ws.Cells.Item(6, 28).Value = 12885.28
$ws.Cells.Item(6, 29).Value = 4813
$ws.Cells.Item(6, 30).Value = 43.53
$ws.Cells.Item(6, 31).Value = 64109
$ws.Cells.Item(6, 32).Value = 3.27
$ws.Cells.Item(6, 33).Value = 1180
$ws.Cells.Item(6, 34).Value = 0.5600000000000001
$ws.Cells.Item(6, 35).Value = 24.52
$ws.Cells.Item(6, 36).Value = 58458490

# Row 7
$ws.Cells.Item(7, 4).Value = 55861
$ws.Cells.Item(7, 5).Value = 4520
$ws.Cells.Item(7, 7).Value = 4711
$ws.Cells.Item(7, 8).Value = 3377
$ws.Cells.Item(7, 9).Value = 3392
$ws.Cells.Item(7, 11).Value = 58138
$ws.Cells.Item(7, 12).Value = 11520
$ws.Cells.Item(7, 13).Value = 46618
$ws.Cells.Item(7, 14).Value = 46506
$ws.Cells.Item(7, 16).Value = 348
$ws.Cells.Item(7, 17).Value = 6627
$ws.Cells.Item(7, 18).Value = -3780
$ws.Cells.Item(7, 19).Value = -1969
$ws.Cells.Item(7, 20).Value = 3069
$ws.Cells.Item(7, 21).Value = 3229
$ws.Cells.Item(7, 23).Value = 8.09
$ws.Cells.Item(7, 24).Value = 6.04
$ws.Cells.Item(7, 25).Value = 7.48
$ws.Cells.Item(7, 26).Value = 6.04
$ws.Cells.Item(7, 27).Value = 24.71
$ws.Cells.Item(7, 29).Value = 4914
$ws.Cells.Item(7, 30).Value = 38.26
$ws.Cells.Item(7, 31).Value = 67891
$ws.Cells.Item(7, 32).Value = 2.77
$ws.Cells.Item(7, 33).Value = 1184
$ws.Cells.Item(7, 34).Value = 0.63
$ws.Cells.Item(7, 35).Value = 20.41

# Row 8
$ws.Cells.Item(8, 4).Value = 60784
$ws.Cells.Item(8, 5).Value = 5783
$ws.Cells.Item(8, 7).Value = 5858
$ws.Cells.Item(8, 8).Value = 4299
$ws.Cells.Item(8, 9).Value = 4293
$ws.Cells.Item(8, 11).Value = 61671
$ws.Cells.Item(8, 12).Value = 11603
$ws.Cells.Item(8, 13).Value = 50068
$ws.Cells.Item(8, 14).Value = 49953
$ws.Cells.Item(8, 16).Value = 348
$ws.Cells.Item(8, 17).Value = 6838
$ws.Cells.Item(8, 18).Value = -3233
$ws.Cells.Item(8, 19).Value = -1106
$ws.Cells.Item(8, 20).Value = 3132
$ws.Cells.Item(8, 21).Value = 3294
$ws.Cells.Item(8, 23).Value = 9.51
$ws.Cells.Item(8, 24).Value = 7.07
$ws.Cells.Item(8, 25).Value = 8.9
$ws.Cells.Item(8, 26).Value = 7.18
$ws.Cells.Item(8, 27).Value = 23.18
$ws.Cells.Item(8, 29).Value = 6220
$ws.Cells.Item(8, 30).Value = 30.23
$ws.Cells.Item(8, 31).Value = 72923
$ws.Cells.Item(8, 32).Value = 2.58
$ws.Cells.Item(8, 33).Value = 1309
$ws.Cells.Item(8, 34).Value = 0.7
$ws.Cells.Item(8, 35).Value = 17.83

# Row 9
$ws.Cells.Item(9, 4).Value = 65462
$ws.Cells.Item(9, 5).Value = 6799
$ws.Cells.Item(9, 7).Value = 6882
$ws.Cells.Item(9, 8).Value = 5066
$ws.Cells.Item(9, 9).Value = 5082
$ws.Cells.Item(9, 11).Value = 65747
$ws.Cells.Item(9, 12).Value = 11497
$ws.Cells.Item(9, 13).Value = 54250
$ws.Cells.Item(9, 14).Value = 54075
$ws.Cells.Item(9, 16).Value = 348
$ws.Cells.Item(9, 17).Value = 7542
$ws.Cells.Item(9, 18).Value = -3801
$ws.Cells.Item(9, 19).Value = -1349
$ws.Cells.Item(9, 20).Value = 3140
$ws.Cells.Item(9, 21).Value = 4160
$ws.Cells.Item(9, 23).Value = 10.39
$ws.Cells.Item(9, 24).Value = 7.74
$ws.Cells.Item(9, 25).Value = 9.77
$ws.Cells.Item(9, 26).Value = 7.95
$ws.Cells.Item(9, 27).Value = 21.19
$ws.Cells.Item(9, 29).Value = 7364
$ws.Cells.Item(9, 30).Value = 25.53
$ws.Cells.Item(9, 31).Value = 78940
$ws.Cells.Item(9, 32).Value = 2.38
$ws.Cells.Item(9, 33).Value = 1417
$ws.Cells.Item(9, 34).Value = 0.75
$ws.Cells.Item(9, 35).Value = 16.29
